$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Merge cells for the "Model" and "Strategy" columns of the new
# rows first (before applying formatting/values), matching the
# merge layout already used for rows 2-9.
$ws.Range("A10:A17").Merge()
$ws.Range("B10:B11").Merge()
$ws.Range("B12:B13").Merge()
$ws.Range("B14:B15").Merge()
$ws.Range("B16:B17").Merge()

# --- Extend the table with 8 new rows (10-17). Columns A-C use the
# same bold/centered/bordered header-row style ("s=1") as the rest of
# the table, so copy that formatting down onto the new rows. Columns
# D-G keep the default (no explicit) style, matching the rest of the
# numeric columns in the table.
$ws.Range("A1").Copy()
$ws.Range("A10:C17").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Update existing values (row 2-9) and populate the new rows (10-17)
$ws.Range("A2").Value = "mistral-nemo_12b"
$ws.Range("B2").Value = "zero_shot"
$ws.Range("C2").Value = "Raw"
$ws.Range("D2").Value = 0.71
$ws.Range("E2").Value = 0.48
$ws.Range("F2").Value = 0.57
$ws.Range("G2").Value = 0.8
$ws.Range("C3").Value = "Role-based"
$ws.Range("D3").Value = 0.61
$ws.Range("E3").Value = 0.65
$ws.Range("F3").Value = 0.63
$ws.Range("G3").Value = 0.78
$ws.Range("B4").Value = "one_shot"
$ws.Range("C4").Value = "Raw"
$ws.Range("D4").Value = 0.53
$ws.Range("E4").Value = 0.64
$ws.Range("F4").Value = 0.58
$ws.Range("G4").Value = 0.74
$ws.Range("C5").Value = "Role-based"
$ws.Range("D5").Value = 0.5600000000000001
$ws.Range("E5").Value = 0.68
$ws.Range("F5").Value = 0.61
$ws.Range("G5").Value = 0.76
$ws.Range("B6").Value = "few_shot"
$ws.Range("C6").Value = "Raw"
$ws.Range("D6").Value = 0.55
$ws.Range("E6").Value = 0.6
$ws.Range("F6").Value = 0.57
$ws.Range("G6").Value = 0.75
$ws.Range("C7").Value = "Role-based"
$ws.Range("D7").Value = 0.5600000000000001
$ws.Range("E7").Value = 0.6
$ws.Range("F7").Value = 0.58
$ws.Range("G7").Value = 0.76
$ws.Range("B8").Value = "auto_cot"
$ws.Range("C8").Value = "Raw"
$ws.Range("D8").Value = 0.71
$ws.Range("E8").Value = 0.44
$ws.Range("F8").Value = 0.55
$ws.Range("G8").Value = 0.79
$ws.Range("C9").Value = "Role-based"
$ws.Range("D9").Value = 0.5600000000000001
$ws.Range("E9").Value = 0.66
$ws.Range("F9").Value = 0.61
$ws.Range("G9").Value = 0.76
$ws.Range("A10").Value = "llama3.2_3b"
$ws.Range("B10").Value = "zero_shot"
$ws.Range("C10").Value = "Raw"
$ws.Range("D10").Value = 0.43
$ws.Range("E10").Value = 0.6899999999999999
$ws.Range("F10").Value = 0.53
$ws.Range("G10").Value = 0.66
$ws.Range("C11").Value = "Role-based"
$ws.Range("D11").Value = 0.54
$ws.Range("E11").Value = 0.67
$ws.Range("F11").Value = 0.6
$ws.Range("G11").Value = 0.75
$ws.Range("B12").Value = "one_shot"
$ws.Range("C12").Value = "Raw"
$ws.Range("D12").Value = 0.5
$ws.Range("E12").Value = 0.66
$ws.Range("F12").Value = 0.57
$ws.Range("G12").Value = 0.72
$ws.Range("C13").Value = "Role-based"
$ws.Range("D13").Value = 0.46
$ws.Range("E13").Value = 0.8
$ws.Range("F13").Value = 0.59
$ws.Range("G13").Value = 0.68
$ws.Range("B14").Value = "few_shot"
$ws.Range("C14").Value = "Raw"
$ws.Range("D14").Value = 0.66
$ws.Range("E14").Value = 0.36
$ws.Range("F14").Value = 0.46
$ws.Range("G14").Value = 0.77
$ws.Range("C15").Value = "Role-based"
$ws.Range("D15").Value = 0.6
$ws.Range("E15").Value = 0.62
$ws.Range("F15").Value = 0.61
$ws.Range("G15").Value = 0.78
$ws.Range("B16").Value = "auto_cot"
$ws.Range("C16").Value = "Raw"
$ws.Range("D16").Value = 0.4
$ws.Range("E16").Value = 0.77
$ws.Range("F16").Value = 0.53
$ws.Range("G16").Value = 0.62
$ws.Range("C17").Value = "Role-based"
$ws.Range("D17").Value = 0.52
$ws.Range("E17").Value = 0.7
$ws.Range("F17").Value = 0.59
$ws.Range("G17").Value = 0.73
